$wb = $excel.ActiveWorkbook

# Sheets: "List of lines" (rId1) and "Instruction" (rId2)
$wsLines = $wb.Worksheets.Item("List of lines")
$wsInstr = $wb.Worksheets.Item("Instruction")

# Add the example MiMIC/CRIMIC # value to the example row on "List of lines"
$wsLines.Range("C3").Value = "MI99999"

# Split the old "Validated" / "Available" status options on "Instruction" into
# four clearer choices: "Available (Validated)" and "Available (Not validated)"
$wsInstr.Range("H3").Value = "Available (Validated)"
$wsInstr.Range("H5").Value = "Available (Not validated)"

# Update selections / active sheet: "Instruction" becomes the active tab,
# with the new last row (H21) selected; "List of lines" keeps C3 selected.
$wsLines.Range("C3").Select() | Out-Null
$wsInstr.Activate() | Out-Null
$wsInstr.Range("H21").Select() | Out-Null
